$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data (rows 2 and 3) for the new
# "Recommender (editor)" / PCI Registered Reports entry.
$ws.Rows("2:3").Insert()

$ws.Range("A2").Value = "Recommender (editor)"
$ws.Range("B2").Value = "Desde 2023"
$ws.Range("C2").Value = "\href{https://rr.peercommunityin.org/}{PCI Registered Reports}"
$ws.Range("D2").Value = "Recommender"
$ws.Range("E2").Value = "Emisión de recomendaciones de las fases 1 y 2"
$ws.Range("E3").Value = "Perfil \href{https://rr.peercommunityin.org/public/user_public_page?userId=1996}{Recommender}"

# Update the "where" (column D) wording for the Guest Editor, Review Editor and
# Journals Incluyen entries, and add a new "Equipo editorial" entry for the
# Review Editor row.
$ws.Range("D4").Value = "Comité Editorial invitado"
$ws.Range("D6").Value = "Equipo editorial"
$ws.Range("D8").Value = "\textit{Journals} Internacionales"

# Row heights (Excel auto-sizes these wrapped-text rows on save).
$ws.Rows("2").RowHeight = 30
$ws.Rows("8").RowHeight = 30

# Leave the selection the way it ended up after the edit.
[void]$ws.Range("D15").Select()
